$wb = $excel.ActiveWorkbook

# --- Sheet "Full results": updated estimates from rerun with all iterations ---
$wsFull = $wb.Worksheets.Item("Full results")

$wsFull.Range("H2").Value = 0.533984259394677
$wsFull.Range("I2").Value = 0.243997347946986
$wsFull.Range("O2").Value = 0.466259480364939

$wsFull.Range("F3").Value = 0.535068033019811
$wsFull.Range("G3").Value = 0.253105913070893

$wsFull.Range("C4").Value = 0.60493272865191
$wsFull.Range("D4").Value = 0.395523933953698
$wsFull.Range("E4").Value = 1.00045666260561
$wsFull.Range("J4").Value = 0.395343395207135
$wsFull.Range("K4").Value = 0.25299038138037
$wsFull.Range("L4").Value = 0.00108327808914099
$wsFull.Range("M4").Value = 0.0709160851578038
$wsFull.Range("N4").Value = 0.254073659469511

# --- Sheet "For plotting": mirror of the Sibcorr / IOLIB / IORAD estimates ---
$wsPlot = $wb.Worksheets.Item("For plotting")

$wsPlot.Range("C2").Value = 0.395343395207135
$wsPlot.Range("D2").Value = 0.327274826249135
$wsPlot.Range("E2").Value = 0.463411964165134

$wsPlot.Range("C3").Value = 0.254073659469511
$wsPlot.Range("D3").Value = 0.182857610814146
$wsPlot.Range("E3").Value = 0.325289708124875

$wsPlot.Range("C4").Value = 0.466259480364939
$wsPlot.Range("D4").Value = 0.39656656692804
$wsPlot.Range("E4").Value = 0.535952393801837
